# "TFs stable, RPT development"
#
# Rename the double-layer CPE frequency-breakpoint parameter from an
# angular-frequency form (wDL / omega_dl, rad/s) to a time-constant form
# (tauDL / tau_dl, s) in both the Negative- and Positive-Electrode
# Parameters sections of the "Parameters" sheet.
#   Negative Electrode Parameters block -> row 36
#   Positive Electrode Parameters block -> row 72
# Both rows have identical layout: B=description, C=code name, D=latex,
# E=value, F=activation energy, G=unit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameters")

# --- Negative Electrode Parameters: row 36 ---
$ws.Range("C36").Value2 = "tauDL"
$ws.Range("D36").Value2 = "\tau_\mathrm{dl}"
$ws.Range("E36").Value2 = 1000000
# Pull the number format (scientific, as used elsewhere in this column)
# from a sibling cell so the style is reused rather than duplicated.
$ws.Range("E37").Copy()
$ws.Range("E36").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("G36").Value2 = "s"

# --- Positive Electrode Parameters: row 72 (mirror of row 36) ---
$ws.Range("C72").Value2 = "tauDL"
$ws.Range("D72").Value2 = "\tau_\mathrm{dl}"
$ws.Range("E72").Value2 = 1000000
$ws.Range("E73").Copy()
$ws.Range("E72").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("G72").Value2 = "s"

# Reflect the user's final selection on the sheet (was E73, now G73).
$ws.Range("G73").Select()
